$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 128.284162
$ws.Range("H2").Value = 384.852486
$ws.Range("I2").Value = 0.02759562102610217
$ws.Range("J2").Value = 0.02759562102610216
$ws.Range("M2").Value = 6.305846
$ws.Range("N2").Value = 18.917538
$ws.Range("O2").Value = 0.01356150511917599
$ws.Range("P2").Value = 0.01356150511917599
$ws.Range("Q2").Value = 808.940169811052
$ws.Range("R2").Value = 7280.461528299468
$ws.Range("S2").Value = 0.0003742381558123251
$ws.Range("T2").Value = 0.0003742381558123249
$ws.Range("G3").Value = 128.284162
$ws.Range("H3").Value = 384.852486
$ws.Range("I3").Value = 0.02759562102610217
$ws.Range("J3").Value = 0.02759562102610216
$ws.Range("O3").Value = 0.392557056479861
$ws.Range("P3").Value = 0.3925570564798609
$ws.Range("Q3").Value = 23415.92390658188
$ws.Range("R3").Value = 210743.3151592368
$ws.Range("S3").Value = 0.01083285576174043
$ws.Range("T3").Value = 0.01083285576174042
$ws.Range("G4").Value = 128.284162
$ws.Range("H4").Value = 384.852486
$ws.Range("I4").Value = 0.02759562102610217
$ws.Range("J4").Value = 0.02759562102610216
$ws.Range("M4").Value = 127.396393
$ws.Range("N4").Value = 382.189179
$ws.Range("O4").Value = 0.2739817680029065
$ws.Range("P4").Value = 0.2739817680029065
$ws.Range("Q4").Value = 16342.93951782767
$ws.Range("R4").Value = 147086.455660449
$ws.Range("S4").Value = 0.007560697037869653
$ws.Range("T4").Value = 0.007560697037869652
$ws.Range("G5").Value = 128.284162
$ws.Range("H5").Value = 384.852486
$ws.Range("I5").Value = 0.02759562102610217
$ws.Range("J5").Value = 0.02759562102610216
$ws.Range("M5").Value = 19.42400133333333
$ws.Range("N5").Value = 58.272004
$ws.Range("O5").Value = 0.04177372766745037
$ws.Range("P5").Value = 0.04177372766745036
$ws.Range("Q5").Value = 2491.791733733549
$ws.Range("R5").Value = 22426.12560360194
$ws.Range("S5").Value = 0.001152771957558559
$ws.Range("T5").Value = 0.001152771957558559
$ws.Range("G6").Value = 128.284162
$ws.Range("H6").Value = 384.852486
$ws.Range("I6").Value = 0.02759562102610217
$ws.Range("J6").Value = 0.02759562102610216
$ws.Range("M6").Value = 129.3233566666667
$ws.Range("N6").Value = 387.97007
$ws.Range("O6").Value = 0.2781259427306063
$ws.Range("P6").Value = 0.2781259427306062
$ws.Range("Q6").Value = 16590.13843701045
$ws.Range("R6").Value = 149311.245933094
$ws.Range("S6").Value = 0.007675058113121205
$ws.Range("T6").Value = 0.007675058113121203
$ws.Range("I7").Value = 0.03577173741430972
$ws.Range("J7").Value = 0.03577173741430972
$ws.Range("M7").Value = 6.305846
$ws.Range("N7").Value = 18.917538
$ws.Range("O7").Value = 0.01356150511917599
$ws.Range("P7").Value = 0.01356150511917599
$ws.Range("Q7").Value = 1048.615478194781
$ws.Range("R7").Value = 9437.539303753027
$ws.Range("S7").Value = 0.0004851186000659805
$ws.Range("T7").Value = 0.0004851186000659804
$ws.Range("I8").Value = 0.03577173741430972
$ws.Range("J8").Value = 0.03577173741430972
$ws.Range("O8").Value = 0.392557056479861
$ws.Range("P8").Value = 0.3925570564798609
$ws.Range("S8").Value = 0.01404244794453194
$ws.Range("T8").Value = 0.01404244794453193
$ws.Range("I9").Value = 0.03577173741430972
$ws.Range("J9").Value = 0.03577173741430972
$ws.Range("M9").Value = 127.396393
$ws.Range("N9").Value = 382.189179
$ws.Range("O9").Value = 0.2739817680029065
$ws.Range("P9").Value = 0.2739817680029065
$ws.Range("Q9").Value = 21185.07644588612
$ws.Range("R9").Value = 190665.6880129751
$ws.Range("S9").Value = 0.009800803861308296
$ws.Range("T9").Value = 0.009800803861308296
$ws.Range("I10").Value = 0.03577173741430972
$ws.Range("J10").Value = 0.03577173741430972
$ws.Range("M10").Value = 19.42400133333333
$ws.Range("N10").Value = 58.272004
$ws.Range("O10").Value = 0.04177372766745037
$ws.Range("P10").Value = 0.04177372766745036
$ws.Range("Q10").Value = 3230.067535206123
$ws.Range("R10").Value = 29070.60781685511
$ws.Range("S10").Value = 0.001494318816936919
$ws.Range("T10").Value = 0.001494318816936919
$ws.Range("I11").Value = 0.03577173741430972
$ws.Range("J11").Value = 0.03577173741430972
$ws.Range("M11").Value = 129.3233566666667
$ws.Range("N11").Value = 387.97007
$ws.Range("O11").Value = 0.2781259427306063
$ws.Range("P11").Value = 0.2781259427306062
$ws.Range("Q11").Value = 21505.51622934827
$ws.Range("R11").Value = 193549.6460641344
$ws.Range("S11").Value = 0.009949048191466589
$ws.Range("T11").Value = 0.009949048191466588
$ws.Range("G12").Value = 1925.076375666667
$ws.Range("H12").Value = 5775.229127000001
$ws.Range("I12").Value = 0.4141094058766164
$ws.Range("J12").Value = 0.4141094058766165
$ws.Range("M12").Value = 6.305846
$ws.Range("N12").Value = 18.917538
$ws.Range("O12").Value = 0.01356150511917599
$ws.Range("P12").Value = 0.01356150511917599
$ws.Range("Q12").Value = 12139.23516319215
$ws.Range("R12").Value = 109253.1164687293
$ws.Range("S12").Value = 0.005615946827694661
$ws.Range("T12").Value = 0.005615946827694661
$ws.Range("G13").Value = 1925.076375666667
$ws.Range("H13").Value = 5775.229127000001
$ws.Range("I13").Value = 0.4141094058766164
$ws.Range("J13").Value = 0.4141094058766165
$ws.Range("O13").Value = 0.392557056479861
$ws.Range("P13").Value = 0.3925570564798609
$ws.Range("Q13").Value = 351387.4294705927
$ws.Range("R13").Value = 3162486.865235334
$ws.Range("S13").Value = 0.1625615694315486
$ws.Range("T13").Value = 0.1625615694315486
$ws.Range("G14").Value = 1925.076375666667
$ws.Range("H14").Value = 5775.229127000001
$ws.Range("I14").Value = 0.4141094058766164
$ws.Range("J14").Value = 0.4141094058766165
$ws.Range("M14").Value = 127.396393
$ws.Range("N14").Value = 382.189179
$ws.Range("O14").Value = 0.2739817680029065
$ws.Range("P14").Value = 0.2739817680029065
$ws.Range("Q14").Value = 245247.7865094463
$ws.Range("R14").Value = 2207230.078585017
$ws.Range("S14").Value = 0.1134584271687086
$ws.Range("T14").Value = 0.1134584271687086
$ws.Range("G15").Value = 1925.076375666667
$ws.Range("H15").Value = 5775.229127000001
$ws.Range("I15").Value = 0.4141094058766164
$ws.Range("J15").Value = 0.4141094058766165
$ws.Range("M15").Value = 19.42400133333333
$ws.Range("N15").Value = 58.272004
$ws.Range("O15").Value = 0.04177372766745037
$ws.Range("P15").Value = 0.04177372766745036
$ws.Range("Q15").Value = 37392.68608771783
$ws.Range("R15").Value = 336534.1747894606
$ws.Range("S15").Value = 0.01729889354561945
$ws.Range("T15").Value = 0.01729889354561944
$ws.Range("G16").Value = 1925.076375666667
$ws.Range("H16").Value = 5775.229127000001
$ws.Range("I16").Value = 0.4141094058766164
$ws.Range("J16").Value = 0.4141094058766165
$ws.Range("M16").Value = 129.3233566666667
$ws.Range("N16").Value = 387.97007
$ws.Range("O16").Value = 0.2781259427306063
$ws.Range("P16").Value = 0.2781259427306062
$ws.Range("Q16").Value = 248957.3387409144
$ws.Range("R16").Value = 2240616.048668229
$ws.Range("S16").Value = 0.1151745689030452
$ws.Range("T16").Value = 0.1151745689030452
$ws.Range("G17").Value = 37.52106933333334
$ws.Range("H17").Value = 112.563208
$ws.Range("I17").Value = 0.008071278587116393
$ws.Range("J17").Value = 0.008071278587116393
$ws.Range("M17").Value = 6.305846
$ws.Range("N17").Value = 18.917538
$ws.Range("O17").Value = 0.01356150511917599
$ws.Range("P17").Value = 0.01356150511917599
$ws.Range("Q17").Value = 236.6020849713227
$ws.Range("R17").Value = 2129.418764741904
$ws.Range("S17").Value = 0.0001094586858774745
$ws.Range("T17").Value = 0.0001094586858774745
$ws.Range("G18").Value = 37.52106933333334
$ws.Range("H18").Value = 112.563208
$ws.Range("I18").Value = 0.008071278587116393
$ws.Range("J18").Value = 0.008071278587116393
$ws.Range("O18").Value = 0.392557056479861
$ws.Range("P18").Value = 0.3925570564798609
$ws.Range("Q18").Value = 6848.783908358976
$ws.Range("R18").Value = 61639.05517523078
$ws.Range("S18").Value = 0.003168437364187343
$ws.Range("T18").Value = 0.003168437364187342
$ws.Range("G19").Value = 37.52106933333334
$ws.Range("H19").Value = 112.563208
$ws.Range("I19").Value = 0.008071278587116393
$ws.Range("J19").Value = 0.008071278587116393
$ws.Range("M19").Value = 127.396393
$ws.Range("N19").Value = 382.189179
$ws.Range("O19").Value = 0.2739817680029065
$ws.Range("P19").Value = 0.2739817680029065
$ws.Range("Q19").Value = 4780.048894569582
$ws.Range("R19").Value = 43020.44005112624
$ws.Range("S19").Value = 0.002211383177342151
$ws.Range("T19").Value = 0.002211383177342151
$ws.Range("G20").Value = 37.52106933333334
$ws.Range("H20").Value = 112.563208
$ws.Range("I20").Value = 0.008071278587116393
$ws.Range("J20").Value = 0.008071278587116393
$ws.Range("M20").Value = 19.42400133333333
$ws.Range("N20").Value = 58.272004
$ws.Range("O20").Value = 0.04177372766745037
$ws.Range("P20").Value = 0.04177372766745036
$ws.Range("Q20").Value = 728.8093007587592
$ws.Range("R20").Value = 6559.283706828833
$ws.Range("S20").Value = 0.0003371673936263238
$ws.Range("T20").Value = 0.0003371673936263237
$ws.Range("G21").Value = 37.52106933333334
$ws.Range("H21").Value = 112.563208
$ws.Range("I21").Value = 0.008071278587116393
$ws.Range("J21").Value = 0.008071278587116393
$ws.Range("M21").Value = 129.3233566666667
$ws.Range("N21").Value = 387.97007
$ws.Range("O21").Value = 0.2781259427306063
$ws.Range("P21").Value = 0.2781259427306062
$ws.Range("Q21").Value = 4852.350631909397
$ws.Range("R21").Value = 43671.15568718456
$ws.Range("S21").Value = 0.002244831966083102
$ws.Range("T21").Value = 0.002244831966083102
$ws.Range("G22").Value = 2391.540242666666
$ws.Range("H22").Value = 7174.620728
$ws.Range("I22").Value = 0.5144519570958551
$ws.Range("J22").Value = 0.5144519570958552
$ws.Range("M22").Value = 6.305846
$ws.Range("N22").Value = 18.917538
$ws.Range("O22").Value = 0.01356150511917599
$ws.Range("P22").Value = 0.01356150511917599
$ws.Range("Q22").Value = 15080.68447305863
$ws.Range("R22").Value = 135726.1602575277
$ws.Range("S22").Value = 0.006976742849725545
$ws.Range("T22").Value = 0.006976742849725545
$ws.Range("G23").Value = 2391.540242666666
$ws.Range("H23").Value = 7174.620728
$ws.Range("I23").Value = 0.5144519570958551
$ws.Range("J23").Value = 0.5144519570958552
$ws.Range("O23").Value = 0.392557056479861
$ws.Range("P23").Value = 0.3925570564798609
$ws.Range("Q23").Value = 436531.8638618149
$ws.Range("R23").Value = 3928786.774756334
$ws.Range("S23").Value = 0.2019517459778526
$ws.Range("T23").Value = 0.2019517459778526
$ws.Range("G24").Value = 2391.540242666666
$ws.Range("H24").Value = 7174.620728
$ws.Range("I24").Value = 0.5144519570958551
$ws.Range("J24").Value = 0.5144519570958552
$ws.Range("M24").Value = 127.396393
$ws.Range("N24").Value = 382.189179
$ws.Range("O24").Value = 0.2739817680029065
$ws.Range("P24").Value = 0.2739817680029065
$ws.Range("Q24").Value = 304673.600630078
$ws.Range("R24").Value = 2742062.405670702
$ws.Range("S24").Value = 0.1409504567576778
$ws.Range("T24").Value = 0.1409504567576778
$ws.Range("G25").Value = 2391.540242666666
$ws.Range("H25").Value = 7174.620728
$ws.Range("I25").Value = 0.5144519570958551
$ws.Range("J25").Value = 0.5144519570958552
$ws.Range("M25").Value = 19.42400133333333
$ws.Range("N25").Value = 58.272004
$ws.Range("O25").Value = 0.04177372766745037
$ws.Range("P25").Value = 0.04177372766745036
$ws.Range("Q25").Value = 46453.28086227765
$ws.Range("R25").Value = 418079.5277604989
$ws.Range("S25").Value = 0.02149057595370911
$ws.Range("T25").Value = 0.02149057595370911
$ws.Range("G26").Value = 2391.540242666666
$ws.Range("H26").Value = 7174.620728
$ws.Range("I26").Value = 0.5144519570958551
$ws.Range("J26").Value = 0.5144519570958552
$ws.Range("M26").Value = 129.3233566666667
$ws.Range("N26").Value = 387.97007
$ws.Range("O26").Value = 0.2781259427306063
$ws.Range("P26").Value = 0.2781259427306062
$ws.Range("Q26").Value = 309282.0117850679
$ws.Range("R26").Value = 2783538.106065611
$ws.Range("S26").Value = 0.1430824355568901
$ws.Range("T26").Value = 0.1430824355568901
